$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Repayment schedule" sheet gets a new (blank) column inserted before
# column N ("Late"), shifting the existing Late/heading/heading/Outstanding
# columns one place to the right (O/P/Q/R) - part of the "Variable
# Instalments" support for RBI loans.
$ws.Columns("N").Insert()

# The freshly inserted column picks up the width of its left neighbour
# (column M), matching the formatting Excel applies to an inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("I18").Select()
